# Update cryptocurrency price (D) and 1h volume change (E) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "26.265.87"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.685.03"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'218.26"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "'0.5239"
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "'0.2707"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").Value = "'0.06413"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("D10").Value = "'21.99"
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("D11").Value = "'0.07502"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").Value = "1.698.92"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").Value = "'4.562"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "'0.5792"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'0.000008447"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "'64.25"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "26.342.95"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "'4.921"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "'188.55"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "'6.186"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "'144.57"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "'7.691"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "'0.1231"
$ws.Range("E26").Value = "  +4.87%  "
$ws.Range("D27").Value = "'15.78"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").Value = "'0.06666"
$ws.Range("E28").Value = "  +13.73%  "
$ws.Range("D29").Value = "'1.348"
$ws.Range("E29").Value = "  +6.53%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").Value = "'3.571"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("D32").Value = "'3.568"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Value = "'1.656"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("D34").Value = "'1.027"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "'0.6197"
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").Value = "'2.400"
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("D37").Value = "'2.697"
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("D38").Value = "'6.379"
$ws.Range("E38").Value = "  +5.74%  "
$ws.Range("D39").Value = "1.104.65"
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").Value = "'0.8771"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "'100.77"
$ws.Range("D44").Value = "1.832.42"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "'0.00000000110"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "'56.66"
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("D47").Value = "'1.007"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'8.142"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "'0.05271"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "'6.040"
$ws.Range("E51").Value = "  +3.09%  "
